# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed figures from the GitHub Actions scrape.
#
# Many Price values look like plain numbers ("314.13", "11.00", "0.0354",
# ...). A bare $ws.Range(...).Value assignment would let Excel's "looks
# like a number" auto-detection store those as numeric cells, which both
# changes their stored type and silently drops significant trailing/
# leading zeros (e.g. "11.00" -> 11, "0.230" -> 0.23). The original
# workbook keeps these as literal text, so for any Price update that
# parses as a number we first force the cell to Text format, assign the
# string, then restore the cell's style to "Normal" so no stray
# number-format is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.233.02"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "2.314.43"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.27%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.982"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.59"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "2.665.02"
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "2.300.72"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "42.104.85"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.60%  "

$ws.Range("E20").Value = "  -1.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "

$ws.Range("E22").Value = "  -6.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("E25").Value = "  -7.96%  "

$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0904"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.09%  "

$ws.Range("E33").Value = "  -5.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.33%  "

$ws.Range("E36").Value = "  -2.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0354"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.230"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "

$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.41%  "

$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.70%  "

